$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order item rows to append (Order_Item_ID, Order_ID, Product_SKU, Quantity, Price, IsDeleted)
$rows = @(
    @(22, 31, "PAMP001",  12, 90,                 0),
    @(23, 31, "PAMR2234", 12, 191.3898305084745,  0),
    @(24, 32, "PAMP001",  10, 75,                 0),
    @(25, 34, "PAMP001",  20, 150,                0),
    @(26, 34, "PAMR4",    19, 55152.39477503628,  0),
    @(27, 34, "PAMP001",  16, 120,                0),
    @(28, 35, "PMME3455", 5,  99.46107784431138,  0),
    @(29, 36, "PAMR4",    2,  5805.515239477503,  0),
    @(30, 36, "CAND234",  6,  5381.165919282511,  0),
    @(31, 37, "PAMP0234", 2,  2974.681818181818,  0),
    @(32, 38, "CAND234",  2,  1793.72197309417,   0),
    @(33, 39, "PAMP0000", 4,  4,                  0),
    @(34, 39, "PAMP003",  7,  126,                0),
    @(35, 39, "PMME3455", 7,  139.2455089820359,  0),
    @(36, 40, "PAMP001",  8,  60,                 0),
    @(37, 40, "PAMP0234", 11, 16360.75,           0),
    @(38, 40, "PAMR2234", 15, 239.2372881355932,  1),
    @(39, 40, "PAMR2234", 1,  15.94915254237288,  1),
    @(40, 40, "PAMR2234", 0,  0,                  0),
    @(41, 41, "CAND234",  2,  1793.72197309417,   0),
    @(42, 41, "PAMP003",  5,  90,                 0),
    @(43, 42, "CAND234",  5,  4484.304932735426,  0),
    @(44, 42, "PAMP003",  5,  90,                 0)
)

$startRow = 23
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}
